# Add a new entry for the space key of the keyboard/encoding table and
# start wiring it up: cell I9 (the "space" slot of the layout grid) gets
# the label "[space]", which introduces a new shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("I9").Value = "[space]"

# Restore the cursor/selection position recorded for this sheet.
$ws.Range("O11").Select()
